$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 33080.668
$ws.Range("J87").Value = 33080.668
$ws.Range("L87").Value = 33080.668
$ws.Range("N87").Value = -35576.668

$ws.Range("H90").Value = 33080.668
$ws.Range("J90").Value = 33080.668
$ws.Range("L90").Value = 99242.00399999999
$ws.Range("N90").Value = -111722.004

$ws.Range("H92").Value = 736.875
$ws.Range("I92").Value = 736.875
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 736.875
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 511.125
$ws.Range("N92").ClearContents()

$ws.Range("H96").Value = 744.8182
$ws.Range("I96").Value = 644
$ws.Range("K96").Value = 1932
$ws.Range("M96").Value = -559

$ws.Range("H99").Value = 271.33334
$ws.Range("I99").Value = 271.33334
$ws.Range("K99").Value = 814.0000200000001
$ws.Range("M99").Value = 683.9999799999999

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2362.0908
$ws.Range("J2").Value = 4896
$ws.Range("L2").Value = 4896
$ws.Range("N2").Value = -5122

$ws.Range("H45").Value = 4866.533
$ws.Range("I45").Value = 4333
$ws.Range("K45").Value = 4333
$ws.Range("M45").Value = -3956

$ws.Range("H74").Value = 964.36365
$ws.Range("I74").Value = 989.4
$ws.Range("K74").Value = 989.4
$ws.Range("M74").Value = -115.4

$ws.Range("H77").Value = 964.36365
$ws.Range("I77").Value = 989.4
$ws.Range("K77").Value = 4947
$ws.Range("M77").Value = -579

$ws.Range("H97").Value = 2577.75
$ws.Range("I97").Value = 2500
$ws.Range("J97").Value = 2603.6667
$ws.Range("K97").Value = 2500
$ws.Range("L97").Value = 2603.6667
$ws.Range("M97").Value = -2004
$ws.Range("N97").Value = -3595.6667

$ws.Range("H116").Value = 2362.0908
$ws.Range("J116").Value = 4896
$ws.Range("L116").Value = 4896
$ws.Range("N116").Value = -9484

$ws.Range("H122").Value = 2002.3334
$ws.Range("I122").Value = 2002.3334
$ws.Range("K122").Value = 6007.0002
$ws.Range("M122").Value = -3557.0002

$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2362.0908
$ws.Range("J3").Value = 4896
$ws.Range("L3").Value = 4896
$ws.Range("N3").Value = -5124

$ws.Range("H20").Value = 1478.6666
$ws.Range("I20").Value = 1291.3334
$ws.Range("J20").Value = 1666
$ws.Range("K20").Value = 1291.3334
$ws.Range("L20").Value = 1666
$ws.Range("M20").Value = -1044.3334
$ws.Range("N20").Value = -2160

$ws.Range("H99").Value = 2006.125
$ws.Range("I99").Value = 1766.5834
$ws.Range("K99").Value = 1766.5834
$ws.Range("M99").Value = -268.5834

$ws.Range("H105").Value = 3121
$ws.Range("I105").Value = 3600
$ws.Range("J105").Value = 2929.4
$ws.Range("K105").Value = 3600
$ws.Range("L105").Value = 2929.4
$ws.Range("M105").Value = -1853
$ws.Range("N105").Value = -6423.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3669.3125
$ws.Range("I31").Value = 2610.4546
$ws.Range("K31").Value = 2610.4546
$ws.Range("M31").Value = -2315.4546

$ws.Range("H34").Value = 3669.3125
$ws.Range("I34").Value = 2610.4546
$ws.Range("K34").Value = 2610.4546
$ws.Range("M34").Value = -2408.4546

$ws.Range("H41").Value = 8642.857

$ws.Range("H50").Value = 21512.834
$ws.Range("J50").Value = 24998.666
$ws.Range("L50").Value = 24998.666
$ws.Range("N50").Value = -26248.666

$ws.Range("H51").Value = 21274.75
$ws.Range("J51").Value = 21274.75
$ws.Range("L51").Value = 21274.75
$ws.Range("N51").Value = -22746.75

$ws.Range("H60").Value = 21666.666
$ws.Range("J60").Value = 21666.666
$ws.Range("L60").Value = 21666.666
$ws.Range("N60").Value = -22688.666

$ws.Range("H61").Value = 21274.75
$ws.Range("J61").Value = 21274.75
$ws.Range("L61").Value = 21274.75
$ws.Range("N61").Value = -21970.75

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 449.625
$ws.Range("I34").Value = 266.66666
$ws.Range("J34").Value = 998.5
$ws.Range("K34").Value = 799.9999799999999
$ws.Range("L34").Value = 2995.5
$ws.Range("M34").Value = -715.9999799999999
$ws.Range("N34").Value = -3163.5

$ws.Range("H39").Value = 2750
$ws.Range("J39").Value = 2750
$ws.Range("L39").Value = 8250
$ws.Range("N39").Value = -8838

$ws.Range("H55").Value = 1800
$ws.Range("I55").Value = 200
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 600
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = -423
$ws.Range("N55").Value = -15354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4853.6665
$ws.Range("I80").Value = 1997.1666
$ws.Range("J80").Value = 10566.667
$ws.Range("K80").Value = 1997.1666
$ws.Range("L80").Value = 10566.667
$ws.Range("M80").Value = -999.1666
$ws.Range("N80").Value = -12562.667

$ws.Range("H83").Value = 4853.6665
$ws.Range("I83").Value = 1997.1666
$ws.Range("J83").Value = 10566.667
$ws.Range("K83").Value = 9985.833000000001
$ws.Range("L83").Value = 52833.335
$ws.Range("M83").Value = -4993.833000000001
$ws.Range("N83").Value = -62817.335

$ws.Range("H97").Value = 2016.6957
$ws.Range("I97").Value = 1929.7778
$ws.Range("J97").Value = 2329.6
$ws.Range("K97").Value = 1929.7778
$ws.Range("L97").Value = 2329.6
$ws.Range("M97").Value = -1433.7778
$ws.Range("N97").Value = -3321.6

$ws.Range("H113").Value = 3785
$ws.Range("I113").Value = 3785
$ws.Range("K113").Value = 3785
$ws.Range("M113").Value = -1615

$ws.Range("H122").Value = 7285.5713
$ws.Range("I122").Value = 6499.8335
$ws.Range("K122").Value = 19499.5005
$ws.Range("M122").Value = -17049.5005

$ws.Range("H126").Value = 7374.25
$ws.Range("I126").Value = 7374.25
$ws.Range("K126").Value = 22122.75
$ws.Range("M126").Value = -19652.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3949.9375
$ws.Range("I46").Value = 3466.6667
$ws.Range("J46").Value = 4571.2856
$ws.Range("K46").Value = 3466.6667
$ws.Range("L46").Value = 4571.2856
$ws.Range("M46").Value = -3278.6667
$ws.Range("N46").Value = -4947.2856

$ws.Range("H61").Value = 212
$ws.Range("I61").Value = 212
$ws.Range("K61").Value = 212
$ws.Range("M61").Value = -10

$ws.Range("H68").Value = 2643.7778
$ws.Range("I68").Value = 1950.25
$ws.Range("J68").Value = 3198.6
$ws.Range("K68").Value = 1950.25
$ws.Range("L68").Value = 3198.6
$ws.Range("M68").Value = -1201.25
$ws.Range("N68").Value = -4696.6

$ws.Range("H71").Value = 2643.7778
$ws.Range("I71").Value = 1950.25
$ws.Range("J71").Value = 3198.6
$ws.Range("K71").Value = 9751.25
$ws.Range("L71").Value = 15993
$ws.Range("M71").Value = -6007.25
$ws.Range("N71").Value = -23481

$ws.Range("H113").Value = 212
$ws.Range("I113").Value = 212
$ws.Range("K113").Value = 212
$ws.Range("M113").Value = 1958

$ws.Range("H122").Value = 6199.8887
$ws.Range("I122").Value = 4799.8335
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 14399.5005
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -11949.5005
$ws.Range("N122").Value = -31900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H75").Value = 24000
$ws.Range("I75").Value = 22000
$ws.Range("K75").Value = 22000
$ws.Range("M75").Value = -21064

$ws.Range("H78").Value = 24000
$ws.Range("I78").Value = 22000
$ws.Range("K78").Value = 66000
$ws.Range("M78").Value = -61320

$ws.Range("H107").Value = 1198.1052
$ws.Range("I107").Value = 1110.1818
$ws.Range("J107").Value = 1319
$ws.Range("K107").Value = 3330.5454
$ws.Range("L107").Value = 3957
$ws.Range("M107").Value = -1410.5454
$ws.Range("N107").Value = -7797

$ws.Range("H113").Value = 815.8889
$ws.Range("I113").Value = 824
$ws.Range("K113").Value = 2472
$ws.Range("M113").Value = -302

$ws.Range("H122").Value = 2845.5386
$ws.Range("J122").Value = 2651
$ws.Range("L122").Value = 7953
$ws.Range("N122").Value = -12853
